# "commit after minor Changes"
#
# The ManageOfferCodePage test-data sheet is updated so the offer-code
# fixture used by the tests changes from "FIRST25"/25 to "FIRST30"/30.
# Making this edit through the Excel UI is what causes Excel to:
#   - move the "active sheet" (and the tabSelected/selection view state)
#     to ManageOfferCodePage,
#   - leave behind a couple of column-width customizations on the
#     (otherwise still empty) C/D columns of that sheet.
# All of that incidental view-state is reproduced below alongside the
# actual value changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ManageOfferCodePage")

$ws.Activate()

$ws.Range("B1").Value = "FIRST30"
$ws.Range("B2").Value = 30

$ws.Range("B2").Select() | Out-Null

$ws.Columns.Item(3).ColumnWidth = 28.666666666666668
$ws.Columns.Item(4).ColumnWidth = 30.333333333333332
